$wb = $excel.ActiveWorkbook

# Map of cell -> new value to apply to both the "展览" and "全部类型" sheets.
$updates = @{
    "G2"  = 65
    "F7"  = 35
    "F8"  = 498
    "F9"  = 39
    "F10" = 1925
    "F11" = 57
    "F13" = 4004
    "F17" = 70
    "F18" = 10
    "F20" = 2724
    "F22" = 372
    "F25" = 54
    "F27" = 45
    "F30" = 34
    "F32" = 144
    "F33" = 1575
    "F34" = 210
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
